$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '97.964.03'
$ws.Cells.Item(2, 5).Value = '  -0.37%  '
$ws.Cells.Item(3, 4).Value = '3.402.22'
$ws.Cells.Item(3, 5).Value = '  -0.63%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '253.51'
$ws.Cells.Item(5, 5).Value = '  -1.17%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '661.57'
$ws.Cells.Item(6, 5).Value = '  -1.16%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.48'
$ws.Cells.Item(7, 5).Value = '  +1.03%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.427'
$ws.Cells.Item(8, 5).Value = '  -2.18%  '
$ws.Cells.Item(9, 2).Value = 'USDC'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.999'
$ws.Cells.Item(9, 5).Value = '  +0.01%  '
$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.04'
$ws.Cells.Item(10, 5).Value = '  -2.78%  '
$ws.Cells.Item(11, 4).Value = '3.397.78'
$ws.Cells.Item(11, 5).Value = '  -0.68%  '
$ws.Cells.Item(12, 2).Value = 'Avalanche'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '44.54'
$ws.Cells.Item(12, 5).Value = '  +5.48%  '
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.209'
$ws.Cells.Item(13, 5).Value = '  -3.37%  '
$ws.Cells.Item(14, 4).Value = '97.666.17'
$ws.Cells.Item(14, 5).Value = '  -0.41%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.12'
$ws.Cells.Item(15, 5).Value = '  -5.02%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000258'
$ws.Cells.Item(16, 5).Value = '  -3.75%  '
$ws.Cells.Item(17, 4).Value = '4.054.90'
$ws.Cells.Item(17, 5).Value = '  -0.05%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '9.19'
$ws.Cells.Item(18, 5).Value = '  +0.89%  '
$ws.Cells.Item(19, 4).Value = '3.410.90'
$ws.Cells.Item(19, 5).Value = '  -0.35%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.21'
$ws.Cells.Item(20, 5).Value = '  +2.42%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.53'
$ws.Cells.Item(21, 5).Value = '  +3.90%  '
$ws.Cells.Item(22, 2).Value = 'Stellar'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.523'
$ws.Cells.Item(22, 5).Value = '  -10.70%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '511.62'
$ws.Cells.Item(23, 5).Value = '  -0.03%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.43'
$ws.Cells.Item(24, 5).Value = '  -1.00%  '
$ws.Cells.Item(25, 5).Value = '  -2.65%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '6.89'
$ws.Cells.Item(26, 5).Value = '  +3.53%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '97.25'
$ws.Cells.Item(27, 5).Value = '  -4.17%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '12.47'
$ws.Cells.Item(28, 5).Value = '  -3.40%  '
$ws.Cells.Item(29, 4).Value = '3.584.43'
$ws.Cells.Item(29, 5).Value = '  -0.48%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '11.98'
$ws.Cells.Item(30, 5).Value = '  +2.58%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.143'
$ws.Cells.Item(31, 5).Value = '  -5.81%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.00'
$ws.Cells.Item(32, 5).Value = '  +0.13%  '
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.72'
$ws.Cells.Item(33, 5).Value = '  +8.03%  '
$ws.Cells.Item(34, 2).Value = 'Cronos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.190'
$ws.Cells.Item(34, 5).Value = '  -3.64%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.998'
$ws.Cells.Item(35, 5).Value = '  -0.25%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.564'
$ws.Cells.Item(36, 5).Value = '  -2.30%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '29.11'
$ws.Cells.Item(37, 5).Value = '  -3.54%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '7.95'
$ws.Cells.Item(38, 5).Value = '  -2.03%  '
$ws.Cells.Item(39, 5).Value = '  -1.81%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '525.22'
$ws.Cells.Item(40, 5).Value = '  -3.24%  '
$ws.Cells.Item(41, 5).Value = '  -1.23%  '
$ws.Cells.Item(42, 5).Value = '  -0.08%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.863'
$ws.Cells.Item(43, 5).Value = '  -2.50%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '24.40'
$ws.Cells.Item(44, 5).Value = '  -1.31%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.74'
$ws.Cells.Item(45, 5).Value = '  -0.06%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0426'
$ws.Cells.Item(46, 5).Value = '  -3.04%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.69'
$ws.Cells.Item(47, 5).Value = '  -3.49%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '5.63'
$ws.Cells.Item(48, 5).Value = '  -4.75%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.61'
$ws.Cells.Item(49, 5).Value = '  -4.60%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.24'
$ws.Cells.Item(50, 5).Value = '  +5.07%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '56.18'
$ws.Cells.Item(51, 5).Value = '  +3.96%  '
